# Generate Report for Archive
#
# The localization-status report was regenerated:
#   1. The "Status" cell that used to read "Ready for handoff" now reads
#      "In Translation" (it shows up on the Overview sheet in both the
#      zh-cn and de-de status columns, and on each language sheet's own
#      Status column).
#   2. Because the new status text is shorter, the "Status" columns were
#      re-sized narrower on all three sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1. Update the status text everywhere it appears -----------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Narrow the "Status" columns to match the regenerated report --------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
